$wb = $excel.ActiveWorkbook

# --- Sheet "error codes": new rows 5-8 + mail?/abort? flipped to "No" ---
$wsErr = $wb.Worksheets.Item("error codes")

$errCodes = @(10201, 10301, 10401, 10501)
$errRow = 5
foreach ($code in $errCodes) {
    $wsErr.Range("A$errRow").Value = $code
    $wsErr.Range("A$errRow").NumberFormat = "000000"
    $wsErr.Range("B$errRow").Value = "Error %SYSTEM%, %SUBSYSTEM%, %LEVEL%"
    $wsErr.Range("D$errRow").Value = "Yes"
    $errRow = $errRow + 1
}

# "mail?" and "abort?" are "No" for the original row 2 and every new row.
$wsErr.Range("C2").Value = "No"
$wsErr.Range("E2").Value = "No"
$wsErr.Range("C5:C8").Value = "No"
$wsErr.Range("E5:E8").Value = "No"

# --- Sheet "database info": new rows 5-8 mirroring row 2 ---
$wsDb = $wb.Worksheets.Item("database info")

$dbCodes = @("010201", "010301", "010401", "010501")
$dbRow = 5
foreach ($code in $dbCodes) {
    $wsDb.Range("A$dbRow").Value = "N"
    $wsDb.Range("B$dbRow").Value = "tb_sales_log"
    $wsDb.Range("C$dbRow").NumberFormat = "@"
    $wsDb.Range("C$dbRow").Value = $code
    $wsDb.Range("D$dbRow").Value = "error_dttm"
    $wsDb.Range("E$dbRow").Value = "error_desc"
    $wsDb.Range("F$dbRow").Value = "transformation_name"
    $dbRow = $dbRow + 1
}

# Update the in-sheet selections to match where the editor left off.
$wsErr.Range("C2").Select()
$wsDb.Range("F5:F8").Select()
$wsDb.Activate()
